# Automatische test-sync: 2025-06-20 16:00:50
# Adds a new incoming mail-log entry (row 18) to the "Logs" sheet,
# extends the conditional formatting ranges to include it, and
# refreshes the "Dashboard" category-count table to reflect the
# updated "Productinformatie" count (now tied with "Openingstijden / Locatie").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Logs sheet: append the new mail entry as row 18
# ---------------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A18").Value = "Is product X op voorraad?"
$logs.Range("B18").Value = "mailmind.test@zohomail.eu"
$logs.Range("C18").Value = "Ik ben geïnteresseerd in product X. Is dit momenteel op voorraad?"
$logs.Range("D18").Value = "Productinformatie"
$logs.Range("F18").Value = "2025-06-20 16:00:16"
$logs.Range("G18").Value = "Nee"

# ---------------------------------------------------------------------------
# 2. Extend the conditional formatting ranges from row 17 to row 18
#    (D2:D17 -> D2:D18, G2:G17 -> G2:G18) while keeping the existing rules,
#    priorities and formatting (dxfId) intact.
# ---------------------------------------------------------------------------
$newRangeD = $logs.Range("D2:D18")
$newRangeG = $logs.Range("G2:G18")

$fcs = $logs.Cells.FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fc = $fcs.Item($i)
    $addr = $fc.AppliesTo.Address()
    if ($addr -eq "`$D`$2:`$D`$17") {
        $fc.ModifyAppliesToRange($newRangeD)
    } elseif ($addr -eq "`$G`$2:`$G`$17") {
        $fc.ModifyAppliesToRange($newRangeG)
    }
}

# ---------------------------------------------------------------------------
# 3. Dashboard sheet: "Productinformatie" now has 2 occurrences (tied with
#    "Openingstijden / Locatie"), so it moves up to row 5 and
#    "Openingstijden / Locatie" shifts down to row 6.
# ---------------------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Range("A5").Value = "Productinformatie"
$dashboard.Range("B5").Value = 2
$dashboard.Range("A6").Value = "Openingstijden / Locatie"
$dashboard.Range("B6").Value = 2
